$wb = $excel.ActiveWorkbook

# Update "想去人数" (F8, F10) on both the "展览" sheet and the "全部类型" sheet.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F8").Value = 136
    $ws.Range("F10").Value = 420
}
